# Applies the "NemoMod baseline" update:
#  - adds three new parameter rows (limit_gnrl_annual_emissions_mt_ch4/n2o/co2)
#    to the "strategy_id-0" sheet, following the existing "General" rows.
#  - updates the sheet's selection/scroll position to reflect where the
#    author ended up after entering the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("strategy_id-0")

# Column letters for the year columns J (year 0) ... AS (year 35), matching
# the layout already used by every other data row on this sheet.
$yearCols = @(
  "J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z",
  "AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN",
  "AO","AP","AQ","AR","AS"
)

$newRows = @(
  @{ Row = 13; Subsector = "General"; Variable = "limit_gnrl_annual_emissions_mt_ch4" },
  @{ Row = 14; Subsector = "General"; Variable = "limit_gnrl_annual_emissions_mt_n2o" },
  @{ Row = 15; Subsector = "General"; Variable = "limit_gnrl_annual_emissions_mt_co2" }
)

foreach ($r in $newRows) {
  $row = $r.Row
  $ws.Cells.Item($row, 1).Value = $r.Subsector
  $ws.Cells.Item($row, 2).Value = $r.Variable
  $ws.Range("H$row").Value = 1
  $ws.Range("I$row").Value = 1
  foreach ($col in $yearCols) {
    $ws.Range("$col$row").Value = -999
  }
}

# Reflect the author's final selection / scroll position on the sheet.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 30
$ws.Range("AQ15").Select()
